# edit.ps1 - apply the "Refine the main document ..." flowchart tweaks
#
# All the touched shapes live inside the single top-level group shape
# ("グループ化 108", id=109) on slide 1. PowerPoint COM automation reports
# shape geometry in points, while the underlying OOXML stores EMUs
# (1 pt = 12700 EMU). The host's Left/Top/Width/Height setters convert the
# point value back to EMU by truncating rather than rounding, so a small
# epsilon is added before conversion to land on the exact target EMU.

function EmuToPt {
    param([double]$Emu)
    return ($Emu / 12700.0) + 0.00004
}

function Get-ShapeById {
    param($Container, [int]$Id)
    for ($i = 1; $i -le $Container.GroupItems.Count; $i++) {
        $candidate = $Container.GroupItems.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    throw "Shape with id $Id not found"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# --- Group container: outer bounding box shrinks by 1 EMU in width ---
$grp.Width = EmuToPt 8672397

# --- Shape id=4: "Experimental Subjects" box ---
$sh4 = Get-ShapeById $grp 4
$sh4.Left = EmuToPt 280800
$sh4.Top = EmuToPt 1107785
$sh4.Width = EmuToPt 5168266
$sh4.Height = EmuToPt 583637
$sh4Para2Run = $sh4.TextFrame.TextRange.Paragraphs(2).Runs(1)
$sh4Para2Run.Text = "The JMDP office sent out an HLA match letter to matched donors."

# --- Shape id=10: connector from id 9 to id 4 ---
$sh10 = Get-ShapeById $grp 10
$sh10.HorizontalFlip = -1
$sh10.Left = EmuToPt 2864933
$sh10.Top = EmuToPt 928751
$sh10.Width = EmuToPt 2
$sh10.Height = EmuToPt 179034

# --- Shape id=16: "CT" box ---
$sh16 = Get-ShapeById $grp 16
$sh16Para1Run = $sh16.TextFrame.TextRange.Paragraphs(1).Runs(1)
$sh16Para1Run.Text = "CT [Primary Outcome] (N = 2,651)"

# --- Shape id=25: "Excluded (N = 105)" box ---
$sh25 = Get-ShapeById $grp 25
$sh25.Left = EmuToPt 5610707
$sh25.Top = EmuToPt 1337007
$sh25.Width = EmuToPt 3317289
$sh25.Height = EmuToPt 1047201

# --- Shape id=33: "Analysis Sample" box ---
$sh33 = Get-ShapeById $grp 33
$sh33.Left = EmuToPt 280799
$sh33.Top = EmuToPt 2083863
$sh33.Width = EmuToPt 5168267
$sh33.Height = EmuToPt 579221
$sh33Para2Run = $sh33.TextFrame.TextRange.Paragraphs(2).Runs(1)
$sh33Para2Run.Text = "Coordination involving a matched donor who lived in Japan."

# --- Shape id=51: connector from id 4 to id 33 ---
$sh51 = Get-ShapeById $grp 51
$sh51.HorizontalFlip = 0
$sh51.Left = EmuToPt 2864933
$sh51.Top = EmuToPt 1691422
$sh51.Width = EmuToPt 0
$sh51.Height = EmuToPt 392441

# --- Shape id=54: connector ending at id 25 ---
$sh54 = Get-ShapeById $grp 54
$sh54.Left = EmuToPt 2864932
$sh54.Top = EmuToPt 1860608
$sh54.Width = EmuToPt 2745775
$sh54.Height = EmuToPt 1472

# --- Shape id=58: connector from id 33 to id 5 ---
$sh58 = Get-ShapeById $grp 58
$sh58.Left = EmuToPt 2864933
$sh58.Top = EmuToPt 2663084
$sh58.Width = EmuToPt 0
$sh58.Height = EmuToPt 611017

# --- Shape id=71: "Dropout (N = 4,907)" box ---
$sh71 = Get-ShapeById $grp 71
$sh71.Left = EmuToPt 5610706
$sh71.Top = EmuToPt 2498663
$sh71.Width = EmuToPt 3342490
$sh71.Height = EmuToPt 1150997

# --- Shape id=93: connector, now anchored to id 71 connection site 1 ---
$sh93 = Get-ShapeById $grp 93
$sh93.ConnectorFormat.EndConnect($sh71, 1)
$sh93.VerticalFlip = 0
$sh93.Left = EmuToPt 2860313
$sh93.Top = EmuToPt 3074162
$sh93.Width = EmuToPt 2750393
$sh93.Height = EmuToPt 0
